# Update davy's status report with new hours entries (hours update src update)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Carry the existing date style (format m/d/yyyy, used by A4:A79) down onto
# the new date cells A80:A85 before filling in the values, so the new rows
# reuse the workbook's existing date style instead of creating a new one.
$ws.Range("A79").Copy()
$ws.Range("A80:A85").PasteSpecial(-4122)

# Row 80
$ws.Range("A80").Value = 40275
$ws.Range("B80").Value = 6
$ws.Range("C80").Value = "beagleboard avr-can interface"

# Row 81
$ws.Range("A81").Value = 40277
$ws.Range("B81").Value = 1
$ws.Range("C81").Value = "buying parts for interface board"

# Row 82
$ws.Range("A82").Value = 40277
$ws.Range("B82").Value = 6.5
$ws.Range("C82").Value = "interface board soldering and assembly"

# Row 83
$ws.Range("A83").Value = 40278
$ws.Range("B83").Value = 8
$ws.Range("C83").Value = "interface board debug"

# Row 84
$ws.Range("A84").Value = 40279
$ws.Range("B84").Value = 3
$ws.Range("C84").Value = "pair programming with Erica"

# Row 85
$ws.Range("A85").Value = 40279
$ws.Range("B85").Value = 6
$ws.Range("C85").Value = "serial interface debug"

# Update the view to match where the author left off editing: scrolled down
# so row 62 is at the top, with the next empty row/cell selected.
$ws.Range("C86").Select()
$excel.ActiveWindow.ScrollRow = 62
$excel.ActiveWindow.ScrollColumn = 1
